$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    4  = "Trust Fund: School"
    21 = "Trust Fund: School"
    23 = "Sovereign Lands, Common Schools, Strategic Investment and Improvement Fund, ND State University, ND Industrial School, School of Mines, Ellendale, School for the Blind, ND School of Science, University of ND, ND State Treasurer, Mayville, Valley City"
    27 = "Penitentiary Land Fund, School Lands"
    28 = "Trust Fund: University, Trust Fund: Swamp, Trust Fund: School"
    29 = "Perm Cmn Schls (Indmty Selec), Perm Common Schls"
    30 = "Perm Common Schls, Agriculture & Mechanical Cllge"
    41 = "Trust Fund: Swamp, Trust Fund: School, Trust Fund: Ind School"
    43 = "Trust Fund: School, Trust Fund: Ind School"
    44 = "Perm Common Schls, Perm Cmn Schls (Indmty Selec), County Bond, New Mexico Institute of Mining and Technology, Common Schools, Public Schools"
    51 = "Trust Fund: School"
    53 = "Trust Fund: Ind School, Trust Fund: Swamp"
    57 = "Perm Cmn Schls (Indmty Selec), 02-90-0021 Fed Patent"
    58 = "Perm Common Schls"
    65 = "Common Schools, Strategic Investment and Improvement Fund, State Hospital, Valley/Mayville, School for the Deaf, ND State University, ND School of Science, Valley City, Veterans Home, Ellendale, Mayville, ND Industrial School"
    67 = "Rural Credit, Indemnity, School and Public Lands, Strategic Investment and Improvement Fund, Common Schools, ND School of Science, Capitol Building, Veterans Home, Mayville, School of Mines, University of ND, ND Industrial School, ND State University, Ellendale, School for the Deaf, State Hospital, Valley/Mayville"
    68 = "County Bond, University, St Chrtbl, Penal & Reform Inst, Perm Common Schls, School of Mines"
    74 = "Trust Fund: Swamp, Trust Fund: School, Trust Fund: Ind School"
    80 = "Perm Cmn Schls (Indmty Selec), Saline Lands, University of New Mexico, Common Schools, New Mexico School for the Visually Handicapped"
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
